$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# 1) Establish formatting for the two brand-new rows (4 and 5) by
#    copying the (border/alignment/number-format) of row 3, which
#    already carries the exact style pattern we need for every new
#    data row: A/B/C/E/F/H = "normal" style, D = centered style,
#    G = right-aligned/number-format style.
# ------------------------------------------------------------------
$ws.Range("A3:H3").Copy()
$ws.Range("A4:H5").PasteSpecial(-4122)   # xlPasteFormats

# ------------------------------------------------------------------
# 2) Row 2 - update the incident that was removed (replaced by a new
#    incident for a different station).
# ------------------------------------------------------------------
$ws.Range("A2").Value = "SR_PTO016M_HNI"
$ws.Range("B2").Value = "Tam-Hiep-Thon-Thuong-PTO_HNI"
$ws.Range("C2").Value = "POWER_AC_EAS"
$ws.Range("D2").Value = "09/05/2025 14:07:23"
$ws.Range("E2").Value = ""
$ws.Range("F2").Value = "Trạm viễn thông loại 2"
$ws.Range("G2").Value = "Phúc Thọ"
$ws.Range("H2").Value = 0.19

# ------------------------------------------------------------------
# 3) Row 3 - update in place as well.
# ------------------------------------------------------------------
$ws.Range("A3").Value = "UL_TTT093M_HNI"
$ws.Range("B3").Value = "THACH-HOA-TTT_HNI"
$ws.Range("C3").Value = "POWER_AC_EAS"
$ws.Range("D3").Value = "09/05/2025 12:39:56"
$ws.Range("E3").Value = ""
$ws.Range("F3").Value = "Trạm viễn thông loại 3"
$ws.Range("G3").Value = "Thạch Thất"
$ws.Range("H3").Value = 1.65

# ------------------------------------------------------------------
# 4) Row 4 - brand-new row.
# ------------------------------------------------------------------
$ws.Range("A4").Value = "4G-STY003M-HNI"
$ws.Range("B4").Value = "Lang-Van-Hoa-STY_HNI"
$ws.Range("C4").Value = "SITE_OOS"
$ws.Range("D4").Value = "08/05/2025 08:56:47"
$ws.Range("E4").Value = "184602- VTHN ĐKTĐ - Thay cột treo anten trạm , dự kiến từ 08h00 ngày 07/05 đến ngày 12/05 - 4 - hanhhh - 08/05/2025 09:17:21"
$ws.Range("F4").Value = "Trạm viễn thông loại 1"
$ws.Range("G4").Value = "Sơn Tây"
$ws.Range("H4").Value = 29.37

# ------------------------------------------------------------------
# 5) Row 5 - brand-new row.
# ------------------------------------------------------------------
$ws.Range("A5").Value = "3G_STY003M_HNI"
$ws.Range("B5").Value = "Lang-Van-Hoa-STY_HNI"
$ws.Range("C5").Value = "SITE_OOS"
$ws.Range("D5").Value = "08/05/2025 08:10:05"
$ws.Range("E5").Value = "184602- VTHN ĐKTĐ - Thay cột treo anten trạm , dự kiến từ 08h00 ngày 07/05 đến ngày 12/05  - 1 - hanhhh - 08/05/2025 08:54:42"
$ws.Range("F5").Value = "Trạm viễn thông loại 1"
$ws.Range("G5").Value = "Sơn Tây"
$ws.Range("H5").Value = 30.15

# ------------------------------------------------------------------
# 6) Column width tweaks (B narrower, E much wider). Excel's COM
#    ColumnWidth setter snaps to a pixel grid, so we pick the input
#    that lands closest to the target stored width (30.7109375 /
#    127.7109375 character-units).
# ------------------------------------------------------------------
$ws.Columns.Item(2).ColumnWidth = 29.8
$ws.Columns.Item(5).ColumnWidth = 126.8
